# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values for data rows 2-10 (row 11 unchanged).
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
